$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.34649999999999
$ws.Range("A8").Value = -22.40720000000002
$ws.Range("A10").Value = -21.61299999999999
$ws.Range("A12").Value = -21.6118
$ws.Range("C12").Value = -10.4596
$ws.Range("C15").Value = -14.6905
$ws.Range("C17").Value = -13.81759999999999
$ws.Range("A18").Value = -22.29200000000001
$ws.Range("C26").Value = -12.38940000000001
$ws.Range("C27").Value = -13.1208
$ws.Range("C28").Value = -13.8309
$ws.Range("A37").Value = -20.23239999999999
$ws.Range("C37").Value = -12.8937
$ws.Range("C47").Value = -12.80819999999999
$ws.Range("A55").Value = -22.15939999999999
$ws.Range("C65").Value = -12.88059999999999
$ws.Range("A68").Value = -21.46739999999999
$ws.Range("C73").Value = -10.96370000000001
$ws.Range("A77").Value = -20.30409999999998
$ws.Range("A78").Value = -19.67309999999998
$ws.Range("A81").Value = -21.99880000000002
$ws.Range("A82").Value = -21.82560000000001
$ws.Range("C84").Value = -12.87589999999999
$ws.Range("C85").Value = -12.3768
$ws.Range("C93").Value = -10.316
$ws.Range("C95").Value = -12.719
$ws.Range("C98").Value = -12.5802
$ws.Range("C99").Value = -11.93870000000001
$ws.Range("C101").Value = -12.9348
